# Apply the edits described by the diff:
# 1. G2 value changes from "S" to "S | R1"
# 2. New row 10 is added with A10 = "Day 9" and B10 = 6/3/2025 (date, formatted like column B)
# 3. Selection changes to C10

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update G2 to reflect a new status value "S | R1"
$ws.Range("G2").Value = "S | R1"

# 2. Add new row 10 data: Day label and date
$ws.Range("A10").Value = "Day 9"
# Copy formatting (date number format) from the cell above so the new
# cell reuses the existing date style instead of creating a new one.
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("B10").Value = (Get-Date -Year 2025 -Month 6 -Day 3).Date

# 3. Update the active selection to C10, matching the diff's sheetView selection
$ws.Range("C10").Select()
